# Apply the content edit described by the commit: add a new row of data
# ("Ravi Ramawat" / "Jaipur") to the "Shedule" worksheet, which is the
# active sheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shedule")

# Make sure this sheet is the active one (it already is tabSelected in the
# source file, but be explicit so the active selection lands here).
$ws.Activate()

# Add the new data row.
$ws.Range("A3").Value = "Ravi Ramawat"
$ws.Range("B3").Value = "Jaipur"

# Move/replace the current selection onto the newly entered cell, matching
# the saved cursor position in the edited workbook.
$ws.Range("B3").Select()
